$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs, Wnt5a, Fzd2, ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Wnt5a"
$ws.Range("C2").Value = "Fzd2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1599003333333333
$ws.Range("H2").Value = 0.479701
$ws.Range("I2").Value = 0.0264777194346773
$ws.Range("J2").Value = 0.02647771943467731
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4403303333333333
$ws.Range("N2").Value = 1.320991
$ws.Range("O2").Value = 0.02587852307343662
$ws.Range("P2").Value = 0.02587852307343662
$ws.Range("Q2").Value = 0.07040896707677778
$ws.Range("R2").Value = 0.633680703691
$ws.Range("S2").Value = 0.0006852042733222778
$ws.Range("T2").Value = 0.0006852042733222779

# Row 3: ECs, Wnt5a, Fzd2, FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Wnt5a"
$ws.Range("C3").Value = "Fzd2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1599003333333333
$ws.Range("H3").Value = 0.479701
$ws.Range("I3").Value = 0.0264777194346773
$ws.Range("J3").Value = 0.02647771943467731
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 15.65098733333333
$ws.Range("N3").Value = 46.952962
$ws.Range("O3").Value = 0.9198195222247485
$ws.Range("P3").Value = 0.9198195222247485
$ws.Range("Q3").Value = 2.502598091595778
$ws.Range("R3").Value = 22.523382824362
$ws.Range("S3").Value = 0.02435472324000582
$ws.Range("T3").Value = 0.02435472324000582

# Row 4: ECs, Wnt5a, Fzd2, M1
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Wnt5a"
$ws.Range("C4").Value = "Fzd2"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1599003333333333
$ws.Range("H4").Value = 0.479701
$ws.Range("I4").Value = 0.0264777194346773
$ws.Range("J4").Value = 0.02647771943467731
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.001510333333333333
$ws.Range("N4").Value = 0.004531
$ws.Range("O4").Value = 0.00008876335118539137
$ws.Range("P4").Value = 0.00008876335118539136
$ws.Range("Q4").Value = 0.0002415028034444445
$ws.Range("R4").Value = 0.002173525231
$ws.Range("S4").Value = 0.000002350251108768524
$ws.Range("T4").Value = 0.000002350251108768524

# Row 5: ECs, Wnt5a, Fzd2, M2
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Wnt5a"
$ws.Range("C5").Value = "Fzd2"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1599003333333333
$ws.Range("H5").Value = 0.479701
$ws.Range("I5").Value = 0.0264777194346773
$ws.Range("J5").Value = 0.02647771943467731
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.008111666666666666
$ws.Range("N5").Value = 0.024335
$ws.Range("O5").Value = 0.0004767283493922972
$ws.Range("P5").Value = 0.0004767283493922972
$ws.Range("Q5").Value = 0.001297058203888889
$ws.Range("R5").Value = 0.011673523835
$ws.Range("S5").Value = 0.00001262267948176606
$ws.Range("T5").Value = 0.00001262267948176606

# Row 6: ECs, Wnt5a, Fzd2, sCs
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Wnt5a"
$ws.Range("C6").Value = "Fzd2"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.1599003333333333
$ws.Range("H6").Value = 0.479701
$ws.Range("I6").Value = 0.0264777194346773
$ws.Range("J6").Value = 0.02647771943467731
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.914341
$ws.Range("N6").Value = 2.743023
$ws.Range("O6").Value = 0.0537364630012372
$ws.Range("P6").Value = 0.0537364630012372
$ws.Range("Q6").Value = 0.1462034306803333
$ws.Range("R6").Value = 1.315830876123
$ws.Range("S6").Value = 0.001422818990758676
$ws.Range("T6").Value = 0.001422818990758676

# Row 7: FAPs, Wnt5a, Fzd2, ECs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Wnt5a"
$ws.Range("C7").Value = "Fzd2"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.879152
$ws.Range("H7").Value = 17.637456
$ws.Range("I7").Value = 0.9735222805653226
$ws.Range("J7").Value = 0.9735222805653228
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.4403303333333333
$ws.Range("N7").Value = 1.320991
$ws.Range("O7").Value = 0.02587852307343662
$ws.Range("P7").Value = 0.02587852307343662
$ws.Range("Q7").Value = 2.588768959877334
$ws.Range("R7").Value = 23.298920638896
$ws.Range("S7").Value = 0.02519331880011434
$ws.Range("T7").Value = 0.02519331880011434

# Row 8: FAPs, Wnt5a, Fzd2, FAPs
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Wnt5a"
$ws.Range("C8").Value = "Fzd2"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 5.879152
$ws.Range("H8").Value = 17.637456
$ws.Range("I8").Value = 0.9735222805653226
$ws.Range("J8").Value = 0.9735222805653228
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 15.65098733333333
$ws.Range("N8").Value = 46.952962
$ws.Range("O8").Value = 0.9198195222247485
$ws.Range("P8").Value = 0.9198195222247485
$ws.Range("Q8").Value = 92.01453348274134
$ws.Range("R8").Value = 828.130801344672
$ws.Range("S8").Value = 0.8954647989847426
$ws.Range("T8").Value = 0.8954647989847427

# Row 9: FAPs, Wnt5a, Fzd2, M1
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Wnt5a"
$ws.Range("C9").Value = "Fzd2"
$ws.Range("D9").Value = "M1"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 5.879152
$ws.Range("H9").Value = 17.637456
$ws.Range("I9").Value = 0.9735222805653226
$ws.Range("J9").Value = 0.9735222805653228
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.001510333333333333
$ws.Range("N9").Value = 0.004531
$ws.Range("O9").Value = 0.00008876335118539137
$ws.Range("P9").Value = 0.00008876335118539136
$ws.Range("Q9").Value = 0.008879479237333334
$ws.Range("R9").Value = 0.079915313136
$ws.Range("S9").Value = 0.00008641310007662284
$ws.Range("T9").Value = 0.00008641310007662284

# Row 10: FAPs, Wnt5a, Fzd2, M2
$ws.Range("A10").Value = "FAPs"
$ws.Range("B10").Value = "Wnt5a"
$ws.Range("C10").Value = "Fzd2"
$ws.Range("D10").Value = "M2"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 5.879152
$ws.Range("H10").Value = 17.637456
$ws.Range("I10").Value = 0.9735222805653226
$ws.Range("J10").Value = 0.9735222805653228
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.008111666666666666
$ws.Range("N10").Value = 0.024335
$ws.Range("O10").Value = 0.0004767283493922972
$ws.Range("P10").Value = 0.0004767283493922972
$ws.Range("Q10").Value = 0.04768972130666667
$ws.Range("R10").Value = 0.42920749176
$ws.Range("S10").Value = 0.0004641056699105311
$ws.Range("T10").Value = 0.0004641056699105312

# Row 11: FAPs, Wnt5a, Fzd2, sCs
$ws.Range("A11").Value = "FAPs"
$ws.Range("B11").Value = "Wnt5a"
$ws.Range("C11").Value = "Fzd2"
$ws.Range("D11").Value = "sCs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 5.879152
$ws.Range("H11").Value = 17.637456
$ws.Range("I11").Value = 0.9735222805653226
$ws.Range("J11").Value = 0.9735222805653228
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.914341
$ws.Range("N11").Value = 2.743023
$ws.Range("O11").Value = 0.0537364630012372
$ws.Range("P11").Value = 0.0537364630012372
$ws.Range("Q11").Value = 5.375549718832001
$ws.Range("R11").Value = 48.379947469488
$ws.Range("S11").Value = 0.05231364401047853
$ws.Range("T11").Value = 0.05231364401047853
